$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Recommandations")
$ws2 = $wb.Worksheets.Item("Top_YTD")

# --- Sheet "Recommandations": update rows 2-41 ---
$ws1.Cells.Item(2,1).Value = "BRVM - SERVICES PUBLICS"
$ws1.Cells.Item(2,2).Value = 0
$ws1.Cells.Item(2,3).Value = 6
$ws1.Cells.Item(2,4).Value = 2416.1
$ws1.Cells.Item(2,5).Value = 100.18
$ws1.Cells.Item(2,6).Value = "🟡 Observer"
$ws1.Cells.Item(2,7).Value = "➖ Neutre"
$ws1.Cells.Item(3,1).Value = "CFAO MOTORS CI"
$ws1.Cells.Item(3,2).Value = 0
$ws1.Cells.Item(3,3).Value = 3
$ws1.Cells.Item(3,4).Value = 1915
$ws1.Cells.Item(3,5).Value = 635
$ws1.Cells.Item(3,6).Value = "🟡 Observer"
$ws1.Cells.Item(3,7).Value = "➖ Neutre"
$ws1.Cells.Item(4,1).Value = "BRVM - AUTRES SECTEURS"
$ws1.Cells.Item(4,2).Value = 0
$ws1.Cells.Item(4,3).Value = 3
$ws1.Cells.Item(4,4).Value = 1911.57
$ws1.Cells.Item(4,5).Value = 637.26
$ws1.Cells.Item(4,6).Value = "🟡 Observer"
$ws1.Cells.Item(4,7).Value = "➖ Neutre"
$ws1.Cells.Item(5,1).Value = "NEI-CEDA CI"
$ws1.Cells.Item(5,2).Value = 0
$ws1.Cells.Item(5,3).Value = 3
$ws1.Cells.Item(5,4).Value = 1790
$ws1.Cells.Item(5,5).Value = 600
$ws1.Cells.Item(5,6).Value = "🟡 Observer"
$ws1.Cells.Item(5,7).Value = "➖ Neutre"
$ws1.Cells.Item(6,1).Value = "UNIWAX CI"
$ws1.Cells.Item(6,2).Value = 0
$ws1.Cells.Item(6,3).Value = 3
$ws1.Cells.Item(6,4).Value = 1710
$ws1.Cells.Item(6,5).Value = 575
$ws1.Cells.Item(6,6).Value = "🟡 Observer"
$ws1.Cells.Item(6,7).Value = "➖ Neutre"
$ws1.Cells.Item(7,1).Value = "AIR LIQUIDE CI"
$ws1.Cells.Item(7,2).Value = 0
$ws1.Cells.Item(7,3).Value = 3
$ws1.Cells.Item(7,4).Value = 1630
$ws1.Cells.Item(7,5).Value = 550
$ws1.Cells.Item(7,6).Value = "🟡 Observer"
$ws1.Cells.Item(7,7).Value = "➖ Neutre"
$ws1.Cells.Item(8,1).Value = "SETAO CI"
$ws1.Cells.Item(8,2).Value = 0
$ws1.Cells.Item(8,3).Value = 3
$ws1.Cells.Item(8,4).Value = 1625
$ws1.Cells.Item(8,5).Value = 535
$ws1.Cells.Item(8,6).Value = "🟡 Observer"
$ws1.Cells.Item(8,7).Value = "➖ Neutre"
$ws1.Cells.Item(9,1).Value = "BRVM - DISTRIBUTION"
$ws1.Cells.Item(9,2).Value = 0
$ws1.Cells.Item(9,3).Value = 3
$ws1.Cells.Item(9,4).Value = 1089.63
$ws1.Cells.Item(9,5).Value = 363.49
$ws1.Cells.Item(9,6).Value = "🟡 Observer"
$ws1.Cells.Item(9,7).Value = "➖ Neutre"
$ws1.Cells.Item(10,1).Value = "BRVM - TRANSPORT"
$ws1.Cells.Item(10,2).Value = 0
$ws1.Cells.Item(10,3).Value = 3
$ws1.Cells.Item(10,4).Value = 1046.41
$ws1.Cells.Item(10,5).Value = 343.91
$ws1.Cells.Item(10,6).Value = "🟡 Observer"
$ws1.Cells.Item(10,7).Value = "➖ Neutre"
$ws1.Cells.Item(11,1).Value = "BRVM - AGRICULTURE"
$ws1.Cells.Item(11,2).Value = 0
$ws1.Cells.Item(11,3).Value = 3
$ws1.Cells.Item(11,4).Value = 961.25
$ws1.Cells.Item(11,5).Value = 320.33
$ws1.Cells.Item(11,6).Value = "🟡 Observer"
$ws1.Cells.Item(11,7).Value = "➖ Neutre"
$ws1.Cells.Item(12,1).Value = "BRVM - INDUSTRIE"
$ws1.Cells.Item(12,2).Value = 0
$ws1.Cells.Item(12,3).Value = 3
$ws1.Cells.Item(12,4).Value = 791.04
$ws1.Cells.Item(12,5).Value = 263.87
$ws1.Cells.Item(12,6).Value = "🟡 Observer"
$ws1.Cells.Item(12,7).Value = "➖ Neutre"
$ws1.Cells.Item(13,1).Value = "BRVM - CONSOMMATION DE BASE"
$ws1.Cells.Item(13,2).Value = 0
$ws1.Cells.Item(13,3).Value = 3
$ws1.Cells.Item(13,4).Value = 652.06
$ws1.Cells.Item(13,5).Value = 217.24
$ws1.Cells.Item(13,6).Value = "🟡 Observer"
$ws1.Cells.Item(13,7).Value = "➖ Neutre"
$ws1.Cells.Item(14,1).Value = "BRVM-PRINCIPAL"
$ws1.Cells.Item(14,2).Value = 0
$ws1.Cells.Item(14,3).Value = 3
$ws1.Cells.Item(14,4).Value = 566.39
$ws1.Cells.Item(14,5).Value = 188.97
$ws1.Cells.Item(14,6).Value = "🟡 Observer"
$ws1.Cells.Item(14,7).Value = "➖ Neutre"
$ws1.Cells.Item(15,1).Value = "BRVM - INDUSTRIELS"
$ws1.Cells.Item(15,2).Value = 0
$ws1.Cells.Item(15,3).Value = 3
$ws1.Cells.Item(15,4).Value = 410.19
$ws1.Cells.Item(15,5).Value = 135.88
$ws1.Cells.Item(15,6).Value = "🟡 Observer"
$ws1.Cells.Item(15,7).Value = "➖ Neutre"
$ws1.Cells.Item(16,1).Value = "BRVM-PRESTIGE"
$ws1.Cells.Item(16,2).Value = 0
$ws1.Cells.Item(16,3).Value = 3
$ws1.Cells.Item(16,4).Value = 385.7
$ws1.Cells.Item(16,5).Value = 129.02
$ws1.Cells.Item(16,6).Value = "🟡 Observer"
$ws1.Cells.Item(16,7).Value = "➖ Neutre"
$ws1.Cells.Item(17,1).Value = "BRVM - FINANCES"
$ws1.Cells.Item(17,2).Value = 0
$ws1.Cells.Item(17,3).Value = 3
$ws1.Cells.Item(17,4).Value = 363.69
$ws1.Cells.Item(17,5).Value = 121.75
$ws1.Cells.Item(17,6).Value = "🟡 Observer"
$ws1.Cells.Item(17,7).Value = "➖ Neutre"
$ws1.Cells.Item(18,1).Value = "BRVM - SERVICES FINANCIERS"
$ws1.Cells.Item(18,2).Value = 0
$ws1.Cells.Item(18,3).Value = 3
$ws1.Cells.Item(18,4).Value = 357.42
$ws1.Cells.Item(18,5).Value = 119.65
$ws1.Cells.Item(18,6).Value = "🟡 Observer"
$ws1.Cells.Item(18,7).Value = "➖ Neutre"
$ws1.Cells.Item(19,1).Value = "BRVM - ENERGIE"
$ws1.Cells.Item(19,2).Value = 0
$ws1.Cells.Item(19,3).Value = 3
$ws1.Cells.Item(19,4).Value = 325.23
$ws1.Cells.Item(19,5).Value = 108.76
$ws1.Cells.Item(19,6).Value = "🟡 Observer"
$ws1.Cells.Item(19,7).Value = "➖ Neutre"
$ws1.Cells.Item(20,1).Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$ws1.Cells.Item(20,2).Value = 0
$ws1.Cells.Item(20,3).Value = 3
$ws1.Cells.Item(20,4).Value = 313.48
$ws1.Cells.Item(20,5).Value = 104.93
$ws1.Cells.Item(20,6).Value = "🟡 Observer"
$ws1.Cells.Item(20,7).Value = "➖ Neutre"
$ws1.Cells.Item(21,1).Value = "BRVM - TELECOMMUNICATIONS"
$ws1.Cells.Item(21,2).Value = 0
$ws1.Cells.Item(21,3).Value = 3
$ws1.Cells.Item(21,4).Value = 276.14
$ws1.Cells.Item(21,5).Value = 92.01
$ws1.Cells.Item(21,6).Value = "🟡 Observer"
$ws1.Cells.Item(21,7).Value = "➖ Neutre"
$ws1.Cells.Item(22,1).Value = "SAFCA CI (SAFC)"
$ws1.Cells.Item(22,2).Value = 3
$ws1.Cells.Item(22,3).Value = 0
$ws1.Cells.Item(22,4).Value = 17.5
$ws1.Cells.Item(22,5).Value = 3.31
$ws1.Cells.Item(22,6).Value = "🟢 Achat"
$ws1.Cells.Item(22,7).Value = "✅ Renforcer"
$ws1.Cells.Item(23,1).Value = "BERNABE CI (BNBC)"
$ws1.Cells.Item(23,2).Value = 2
$ws1.Cells.Item(23,3).Value = 0
$ws1.Cells.Item(23,4).Value = 9.33
$ws1.Cells.Item(23,5).Value = 7.14
$ws1.Cells.Item(23,6).Value = "🟡 Observer"
$ws1.Cells.Item(23,7).Value = "➖ Neutre"
$ws1.Cells.Item(24,1).Value = "ORAGROUP TOGO (ORGT)"
$ws1.Cells.Item(24,2).Value = 2
$ws1.Cells.Item(24,3).Value = 0
$ws1.Cells.Item(24,4).Value = 6.23
$ws1.Cells.Item(24,5).Value = 3.7
$ws1.Cells.Item(24,6).Value = "🟡 Observer"
$ws1.Cells.Item(24,7).Value = "➖ Neutre"
$ws1.Cells.Item(25,1).Value = "VIVO ENERGY CI (SHEC)"
$ws1.Cells.Item(25,2).Value = 1
$ws1.Cells.Item(25,3).Value = 0
$ws1.Cells.Item(25,4).Value = 5.09
$ws1.Cells.Item(25,5).Value = 5.09
$ws1.Cells.Item(25,6).Value = "🟡 Observer"
$ws1.Cells.Item(25,7).Value = "➖ Neutre"
$ws1.Cells.Item(26,1).Value = "PALM CI (PALC)"
$ws1.Cells.Item(26,2).Value = 1
$ws1.Cells.Item(26,3).Value = 0
$ws1.Cells.Item(26,4).Value = 3.41
$ws1.Cells.Item(26,5).Value = 3.41
$ws1.Cells.Item(26,6).Value = "🟡 Observer"
$ws1.Cells.Item(26,7).Value = "➖ Neutre"
$ws1.Cells.Item(27,1).Value = "SITAB CI (STBC)"
$ws1.Cells.Item(27,2).Value = 1
$ws1.Cells.Item(27,3).Value = 0
$ws1.Cells.Item(27,4).Value = 3.06
$ws1.Cells.Item(27,5).Value = 3.06
$ws1.Cells.Item(27,6).Value = "🟡 Observer"
$ws1.Cells.Item(27,7).Value = "➖ Neutre"
$ws1.Cells.Item(28,1).Value = "SERVAIR ABIDJAN CI (ABJC)"
$ws1.Cells.Item(28,2).Value = 1
$ws1.Cells.Item(28,3).Value = 1
$ws1.Cells.Item(28,4).Value = 2.65
$ws1.Cells.Item(28,5).Value = -1.42
$ws1.Cells.Item(28,6).Value = "🟡 Observer"
$ws1.Cells.Item(28,7).Value = "👀 À surveiller"
$ws1.Cells.Item(29,1).Value = "ECOBANK TRANS. INCORP. TG (ETIT)"
$ws1.Cells.Item(29,2).Value = 1
$ws1.Cells.Item(29,3).Value = 1
$ws1.Cells.Item(29,4).Value = 0.42
$ws1.Cells.Item(29,5).Value = 6.67
$ws1.Cells.Item(29,6).Value = "🟡 Observer"
$ws1.Cells.Item(29,7).Value = "👀 À surveiller"
$ws1.Cells.Item(30,1).Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$ws1.Cells.Item(30,2).Value = 1
$ws1.Cells.Item(30,3).Value = 1
$ws1.Cells.Item(30,4).Value = 0.33
$ws1.Cells.Item(30,5).Value = 4.04
$ws1.Cells.Item(30,6).Value = "🟡 Observer"
$ws1.Cells.Item(30,7).Value = "👀 À surveiller"
$ws1.Cells.Item(31,1).Value = "UNIWAX CI (UNXC)"
$ws1.Cells.Item(31,2).Value = 1
$ws1.Cells.Item(31,3).Value = 1
$ws1.Cells.Item(31,4).Value = 0.28
$ws1.Cells.Item(31,5).Value = 5.45
$ws1.Cells.Item(31,6).Value = "🟡 Observer"
$ws1.Cells.Item(31,7).Value = "👀 À surveiller"
$ws1.Cells.Item(32,1).Value = "TOTAL"
$ws1.Cells.Item(32,2).Value = 0
$ws1.Cells.Item(32,3).Value = 3
$ws1.Cells.Item(32,4).Value = 0
$ws1.Cells.Item(32,5).Value = 0
$ws1.Cells.Item(32,6).Value = "🟡 Observer"
$ws1.Cells.Item(32,7).Value = "➖ Neutre"
$ws1.Cells.Item(33,1).Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$ws1.Cells.Item(33,2).Value = 0
$ws1.Cells.Item(33,3).Value = 1
$ws1.Cells.Item(33,4).Value = -1.75
$ws1.Cells.Item(33,5).Value = -1.75
$ws1.Cells.Item(33,6).Value = "🟡 Observer"
$ws1.Cells.Item(33,7).Value = "➖ Neutre"
$ws1.Cells.Item(34,1).Value = "FILTISAC CI (FTSC)"
$ws1.Cells.Item(34,2).Value = 0
$ws1.Cells.Item(34,3).Value = 1
$ws1.Cells.Item(34,4).Value = -1.84
$ws1.Cells.Item(34,5).Value = -1.84
$ws1.Cells.Item(34,6).Value = "🟡 Observer"
$ws1.Cells.Item(34,7).Value = "➖ Neutre"
$ws1.Cells.Item(35,1).Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$ws1.Cells.Item(35,2).Value = 0
$ws1.Cells.Item(35,3).Value = 1
$ws1.Cells.Item(35,4).Value = -2
$ws1.Cells.Item(35,5).Value = -2
$ws1.Cells.Item(35,6).Value = "🟡 Observer"
$ws1.Cells.Item(35,7).Value = "➖ Neutre"
$ws1.Cells.Item(36,1).Value = "CIE CI (CIEC)"
$ws1.Cells.Item(36,2).Value = 0
$ws1.Cells.Item(36,3).Value = 1
$ws1.Cells.Item(36,4).Value = -2.27
$ws1.Cells.Item(36,5).Value = -2.27
$ws1.Cells.Item(36,6).Value = "🟡 Observer"
$ws1.Cells.Item(36,7).Value = "➖ Neutre"
$ws1.Cells.Item(37,1).Value = "BANK OF AFRICA NG (BOAN)"
$ws1.Cells.Item(37,2).Value = 0
$ws1.Cells.Item(37,3).Value = 1
$ws1.Cells.Item(37,4).Value = -2.71
$ws1.Cells.Item(37,5).Value = -2.71
$ws1.Cells.Item(37,6).Value = "🟡 Observer"
$ws1.Cells.Item(37,7).Value = "➖ Neutre"
$ws1.Cells.Item(38,1).Value = "SOGB CI (SOGC)"
$ws1.Cells.Item(38,2).Value = 0
$ws1.Cells.Item(38,3).Value = 1
$ws1.Cells.Item(38,4).Value = -3.23
$ws1.Cells.Item(38,5).Value = -3.23
$ws1.Cells.Item(38,6).Value = "🟡 Observer"
$ws1.Cells.Item(38,7).Value = "➖ Neutre"
$ws1.Cells.Item(39,1).Value = "TOTALENERGIES MARKETING CI (TTLC)"
$ws1.Cells.Item(39,2).Value = 1
$ws1.Cells.Item(39,3).Value = 1
$ws1.Cells.Item(39,4).Value = -3.46
$ws1.Cells.Item(39,5).Value = 3.09
$ws1.Cells.Item(39,6).Value = "🟡 Observer"
$ws1.Cells.Item(39,7).Value = "👀 À surveiller"
$ws1.Cells.Item(40,1).Value = "SOLIBRA CI (SLBC)"
$ws1.Cells.Item(40,2).Value = 0
$ws1.Cells.Item(40,3).Value = 2
$ws1.Cells.Item(40,4).Value = -5.2
$ws1.Cells.Item(40,5).Value = -1.97
$ws1.Cells.Item(40,6).Value = "🟡 Observer"
$ws1.Cells.Item(40,7).Value = "➖ Neutre"
$ws1.Cells.Item(41,1).Value = "SAPH CI (SPHC)"
$ws1.Cells.Item(41,2).Value = 0
$ws1.Cells.Item(41,3).Value = 2
$ws1.Cells.Item(41,4).Value = -5.45
$ws1.Cells.Item(41,5).Value = -2.69
$ws1.Cells.Item(41,6).Value = "🟡 Observer"
$ws1.Cells.Item(41,7).Value = "➖ Neutre"

# Remove rows 42-44 (no longer present after the update)
$ws1.Rows.Item(42).Delete()
$ws1.Rows.Item(42).Delete()
$ws1.Rows.Item(42).Delete()

# --- Sheet "Top_YTD": update rows 2-11 ---
$ws2.Cells.Item(2,1).Value = "BRVM - SERVICES PUBLICS"
$ws2.Cells.Item(2,2).Value = 418442.59
$ws2.Cells.Item(3,1).Value = "CFAO MOTORS CI"
$ws2.Cells.Item(3,2).Value = 40141.25
$ws2.Cells.Item(4,1).Value = "BRVM - AUTRES SECTEURS"
$ws2.Cells.Item(4,2).Value = 39961.59
$ws2.Cells.Item(5,1).Value = "NEI-CEDA CI"
$ws2.Cells.Item(5,2).Value = 33711.75
$ws2.Cells.Item(6,1).Value = "UNIWAX CI"
$ws2.Cells.Item(6,2).Value = 29954.38
$ws2.Cells.Item(7,1).Value = "AIR LIQUIDE CI"
$ws2.Cells.Item(7,2).Value = 26524
$ws2.Cells.Item(8,1).Value = "SETAO CI"
$ws2.Cells.Item(8,2).Value = 26303.3
$ws2.Cells.Item(9,1).Value = "BRVM - DISTRIBUTION"
$ws2.Cells.Item(9,2).Value = 9837.69
$ws2.Cells.Item(10,1).Value = "BRVM - TRANSPORT"
$ws2.Cells.Item(10,2).Value = 8939.12
$ws2.Cells.Item(11,1).Value = "BRVM - AGRICULTURE"
$ws2.Cells.Item(11,2).Value = 7330.49
